$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D2"  = "R5 (60)"
    "F2"  = "Tuesday 14:30-15:55, Friday 16:00-17:25 (R)"
    "D3"  = "R7 (60)"
    "F3"  = "Monday 16:00-17:25, Thursday 14:30-15:55 (Q)"
    "D4"  = "R5 (60)"
    "D5"  = "R2 (60)"
    "F5"  = "Monday 11:00-11:55, Wednesday 10:00-10:55, Thursday 9:00-9:55 (C)"
    "D6"  = "R5 (60)"
    "F6"  = "Monday 14:30-15:55, Thursday 16:00-17:25 (P)"
    "D7"  = "R1 (60)"
    "F7"  = "Monday 9:00-9:55, Wednesday 11:00-11:55, Thursday 10:00-10:55 (A)"
    "D8"  = "R4 (60)"
    "F8"  = "Monday 9:00-9:55, Wednesday 11:00-11:55, Thursday 10:00-10:55 (A)"
    "D9"  = "R1 (60)"
    "F9"  = "Tuesday 10:00-10:55, Thursday 11:00-11:55, Friday 10:00-10:55 (F)"
    "D10" = "R2 (60)"
    "F10" = "Monday 12:00-12:55, Tuesday 9:00-9:55, Friday 11:00-11:55 (D)"
    "D11" = "R2 (60)"
    "F11" = "Monday 14:30-15:55, Thursday 16:00-17:25 (P)"
    "D13" = "R5 (60)"
    "F13" = "Tuesday 16:00-17:25, Friday 14:30-15:55 (S)"
    "D14" = "R1 (60)"
    "F14" = "Tuesday 14:30-15:55, Friday 16:00-17:25 (R)"
    "F15" = "Monday 9:00-9:55, Wednesday 11:00-11:55, Thursday 10:00-10:55 (A)"
    "D16" = "R4 (60)"
    "D17" = "R7 (60)"
    "F17" = "Monday 10:00-10:55, Wednesday 9:00-9:55, Friday 9:00-9:55 (B)"
    "D18" = "R3 (60)"
    "F18" = "Tuesday 10:00-10:55, Thursday 11:00-11:55, Friday 10:00-10:55 (F)"
    "F19" = "Monday 12:00-12:55, Tuesday 9:00-9:55, Friday 11:00-11:55 (D)"
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
